$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, border, center alignment) from the existing
# header cell H1 onto the two new header cells I1 and J1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Values for column I (I0) and column J (IF), rows 2-26
$dataI = @(6, 6, 7, 9, 7, 2, 5, 7, 8, 8, 7, 9, 6, 6, 9, 6, 6, 7, 7, 9, 7, 5, 5, 5, 4)
$dataJ = @(6, 6, 7, 7, 7, 3, 5, 8, 8, 8, 7, 9, 6, 6, 9, 6, 7, 8, 8, 9, 7, 6, 6, 5, 4)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
